# Apply the "Deploy: Update site content" edit to the ValueSet workbook.
#
# Changes:
#   1. Metadata sheet: bump the "Date" value.
#   2. Include #0 sheet: replace the three old SNOMED concept rows
#      (370996005/385634002/118222006) with five new concept rows
#      (385652002/385651009/255609007/723510000/260388008), keep the
#      blank separator row, and move the "System URI" row down below
#      the new concept rows.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Metadata")
$ws2 = $wb.Worksheets.Item("Include #0")

# --- 1. Metadata!B8 ("Date" row) -----------------------------------------
$ws1.Range("B8").Value = "2026-01-01T13:37:23+00:00"

# --- 2. Include #0 concept table ------------------------------------------
# Force column A/B on the data rows to keep text formatting, so that the
# all-digit SNOMED codes are stored as text (matching the existing
# "code as text" convention used throughout this sheet) instead of being
# auto-coerced to numbers.
$ws2.Range("A2:B8").NumberFormat = "@"

$concepts = @(
    @("385652002", "Objective achieved"),
    @("385651009", "Objective not achieved"),
    @("255609007", "Partial achievement"),
    @("723510000", "Sustained improvement"),
    @("260388008", "Worsening")
)

$row = 2
foreach ($pair in $concepts) {
    $ws2.Cells.Item($row, 1).Value = $pair[0]
    $ws2.Cells.Item($row, 2).Value = $pair[1]
    $row = $row + 1
}

# Row 7 stays blank (separator row), matching the sheet's existing layout.
$ws2.Cells.Item(7, 1).Value = ""
$ws2.Cells.Item(7, 2).Value = ""

# Row 8: "System URI" / "http://snomed.info/sct" moved down from row 5/6.
$ws2.Cells.Item(8, 1).Value = "System URI"
$ws2.Cells.Item(8, 2).Value = "http://snomed.info/sct"
